$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colA = @(
    -0.07563215593988559,
    -0.08321336037471916,
    -0.09318794296272337,
    -0.0849933915065666,
    -0.08162374674555117,
    0.018860305832472335,
    0.02894504041560353,
    0.0389645132447507,
    0.041038342977747,
    0.043120849825807994,
    0.0009945509405886455,
    0.004494086662748309,
    0.008000319910962617,
    0.01601128148459985,
    0.017041805791581766,
    0.019087407532472955,
    0.021290105672449933,
    0.006669596945126699,
    -0.012090156322797974,
    -0.008015179525065363,
    -0.004005403346071468,
    0.014613891461907613,
    0.019755483766768833,
    -0.02009743398695729,
    -0.09723277202765424,
    -0.09461181515333372,
    -0.09195400105284035,
    -0.08900576131793603,
    -0.08135690326955647,
    -0.02116980688056458,
    -0.014020452618883184,
    -0.004000438144025509
)
$colB = @(
    0.07522860108178264,
    0.08257405143158447,
    0.09299339137253959,
    0.08462374667357508,
    0.08037246611573767,
    -0.018945040606917374,
    -0.028964513430980166,
    -0.03903834304009468,
    -0.04112084987623632,
    -0.043125737608711745,
    -0.0009940867340647586,
    -0.004500319980141221,
    -0.008011281616040478,
    -0.01604180582191539,
    -0.017087407573447955,
    -0.019290105704546257,
    -0.021377727898507892,
    -0.006789275425433061,
    0.012015179461842163,
    0.008005403282041357,
    0.003999999935691001,
    -0.014755483848691853,
    -0.0200284492189029,
    0.019999999706506344,
    0.0971118150617496,
    0.09445400095402334,
    0.09100576119300552,
    0.08835690305161847,
    0.08116980590382461,
    0.021020452382436883,
    0.014000437864311266,
    0.003999999804859655
)

for ($i = 0; $i -lt $colA.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $colA[$i]
    $ws.Cells.Item($row, 2).Value = $colB[$i]
}

# Column widths (closest achievable match to target stored widths
# 15.7109375 / 16.42578125 given the engine's pixel-quantized model)
$ws.Columns.Item(1).ColumnWidth = 14.833333333333334
$ws.Columns.Item(2).ColumnWidth = 15.666666666666666